$d = $word.ActiveDocument

# Locate the paragraph that contains the astromap credit/link line.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*CzechGlobe*GaNight*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range

    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' +
        '<w:p w14:paraId="39E082B4" w14:textId="2B18FFC7" w:rsidR="004615A9" w:rsidRPr="00DB0F3B" w:rsidRDefault="004615A9" w:rsidP="00FA27C1">' +
        '<w:pPr><w:pStyle w:val="BasicParagraph"/><w:pBdr><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="10" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:pBdr><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:right="-90"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="CastleT-Book"/><w:sz w:val="19"/></w:rPr></w:pPr>' +
        '<w:r/><w:r><w:t>Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($xmlFrag)
}
